# Apply BOM updates: add C15 to the C3,C4,C5,C6,C13,C14 capacitor group (row 4)
# and add R12 to the R5,R8,R9,R10 resistor group (row 16), then append a bold
# grand-total row (row 25) summing the Total column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: capacitor group C3,C4,C5,C6,C13,C14 -> add C15, qty 5 -> 7
$ws.Cells.Item(4, 2).Value = "C3, C4, C5, C6, C13, C14, C15"
$ws.Cells.Item(4, 7).Value = 7

# Row 16: resistor group R5,R8,R9,R10 -> add R12, qty 4 -> 5
$ws.Cells.Item(16, 2).Value = "R5, R8, R9, R10, R12"
$ws.Cells.Item(16, 7).Value = 5

# New row 25: grand total, bold
$ws.Cells.Item(25, 9).Formula = "=SUM(I3:I24)"
$ws.Cells.Item(25, 9).Font.Bold = $true

# Restore selection similar to the saved workbook state
$ws.Range("K21").Select()
